$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to text so numeric-looking IDs are stored as strings (matches original inlineStr data)
$ws.Range("A2:A283").NumberFormat = "@"

$arrA = New-Object "object[,]" 282,1
$arrB = New-Object "object[,]" 282,1
$arrC = New-Object "object[,]" 282,1

$arrA[0,0] = "298"
$arrB[0,0] = 0.0018
$arrC[0,0] = 0.3676
$arrA[1,0] = "297"
$arrB[1,0] = -0.0302
$arrC[1,0] = 0.6565
$arrA[2,0] = "296"
$arrB[2,0] = -0.0643
$arrC[2,0] = 0.6147
$arrA[3,0] = "295"
$arrB[3,0] = -0.0163
$arrC[3,0] = 0.0958
$arrA[4,0] = "293"
$arrB[4,0] = -0.0191
$arrC[4,0] = 0.1281
$arrA[5,0] = "292"
$arrB[5,0] = -0.0002
$arrC[5,0] = 0.1981
$arrA[6,0] = "291"
$arrB[6,0] = 0.0015
$arrC[6,0] = 0.3166
$arrA[7,0] = "290"
$arrB[7,0] = -0.0038
$arrC[7,0] = 0.742
$arrA[8,0] = "288"
$arrB[8,0] = 0.0179
$arrC[8,0] = 0.4822
$arrA[9,0] = "287"
$arrB[9,0] = 0.0316
$arrC[9,0] = 0.6574
$arrA[10,0] = "284"
$arrB[10,0] = -0.0112
$arrC[10,0] = 0.1015
$arrA[11,0] = "281"
$arrB[11,0] = 0.0186
$arrC[11,0] = 0.1738
$arrA[12,0] = "280"
$arrB[12,0] = -0.0844
$arrC[12,0] = 0.1891
$arrA[13,0] = "279"
$arrB[13,0] = 0.0557
$arrC[13,0] = 0.2566
$arrA[14,0] = "278"
$arrB[14,0] = 0.0015
$arrC[14,0] = 0.3175
$arrA[15,0] = "276"
$arrB[15,0] = 0.0031
$arrC[15,0] = 0.1311
$arrA[16,0] = "275"
$arrB[16,0] = 0.0831
$arrC[16,0] = 0.1619
$arrA[17,0] = "274"
$arrB[17,0] = -0.0156
$arrC[17,0] = 0.2041
$arrA[18,0] = "271"
$arrB[18,0] = 0.0515
$arrC[18,0] = 0.1889
$arrA[19,0] = "267"
$arrB[19,0] = 0.0025
$arrC[19,0] = 0.2421
$arrA[20,0] = "266"
$arrB[20,0] = -0.0304
$arrC[20,0] = 0.4378
$arrA[21,0] = "264"
$arrB[21,0] = -0.0044
$arrC[21,0] = 0.2959
$arrA[22,0] = "263"
$arrB[22,0] = 0.0709
$arrC[22,0] = 0.2036
$arrA[23,0] = "262"
$arrB[23,0] = 0.0131
$arrC[23,0] = 0.1006
$arrA[24,0] = "261"
$arrB[24,0] = 0.1376
$arrC[24,0] = 0.3925
$arrA[25,0] = "260"
$arrB[25,0] = -0.0116
$arrC[25,0] = 0.3037
$arrA[26,0] = "258"
$arrB[26,0] = -0.0415
$arrC[26,0] = 0.5963
$arrA[27,0] = "257"
$arrB[27,0] = -0.0157
$arrC[27,0] = 0.4804
$arrA[28,0] = "256"
$arrB[28,0] = -0.0242
$arrC[28,0] = 0.1822
$arrA[29,0] = "255"
$arrB[29,0] = 0.027
$arrC[29,0] = 0.1701
$arrA[30,0] = "254"
$arrB[30,0] = 0.2372
$arrC[30,0] = 0.56
$arrA[31,0] = "253"
$arrB[31,0] = 0.0616
$arrC[31,0] = 0.2351
$arrA[32,0] = "252"
$arrB[32,0] = 0.0852
$arrC[32,0] = 0.2099
$arrA[33,0] = "251"
$arrB[33,0] = 0.0238
$arrC[33,0] = 0.2664
$arrA[34,0] = "8"
$arrB[34,0] = 0.0719
$arrC[34,0] = 0.1947
$arrA[35,0] = "6"
$arrB[35,0] = 0.0322
$arrC[35,0] = 0.3508
$arrA[36,0] = "5"
$arrB[36,0] = -0.0208
$arrC[36,0] = 0.1824
$arrA[37,0] = "49"
$arrB[37,0] = -0.0372
$arrC[37,0] = 0.4311
$arrA[38,0] = "48"
$arrB[38,0] = -0.0228
$arrC[38,0] = 0.0863
$arrA[39,0] = "47"
$arrB[39,0] = 0.0312
$arrC[39,0] = 0.2015
$arrA[40,0] = "46"
$arrB[40,0] = -0.0612
$arrC[40,0] = 0.1231
$arrA[41,0] = "45"
$arrB[41,0] = 0.0935
$arrC[41,0] = 1.6582
$arrA[42,0] = "43"
$arrB[42,0] = -0.0052
$arrC[42,0] = 0.2578
$arrA[43,0] = "41"
$arrB[43,0] = 0.0343
$arrC[43,0] = 0.3642
$arrA[44,0] = "40"
$arrB[44,0] = 0.0295
$arrC[44,0] = 0.1025
$arrA[45,0] = "4"
$arrB[45,0] = -0.0237
$arrC[45,0] = 0.1463
$arrA[46,0] = "39"
$arrB[46,0] = 0.0556
$arrC[46,0] = 0.1401
$arrA[47,0] = "38"
$arrB[47,0] = -0.0055
$arrC[47,0] = 0.2114
$arrA[48,0] = "37"
$arrB[48,0] = -0.0091
$arrC[48,0] = 0.5024
$arrA[49,0] = "34"
$arrB[49,0] = 0.0016
$arrC[49,0] = 0.4061
$arrA[50,0] = "33"
$arrB[50,0] = 0.0819
$arrC[50,0] = 0.2125
$arrA[51,0] = "32"
$arrB[51,0] = 0.0132
$arrC[51,0] = 0.2106
$arrA[52,0] = "30"
$arrB[52,0] = 0.0335
$arrC[52,0] = 0.2124
$arrA[53,0] = "3"
$arrB[53,0] = -0.1532
$arrC[53,0] = 0.5537
$arrA[54,0] = "29"
$arrB[54,0] = -0.0161
$arrC[54,0] = 0.2852
$arrA[55,0] = "25"
$arrB[55,0] = 0.0312
$arrC[55,0] = 0.0812
$arrA[56,0] = "23"
$arrB[56,0] = 0.0219
$arrC[56,0] = 0.3161
$arrA[57,0] = "22"
$arrB[57,0] = 0.1666
$arrC[57,0] = 0.3027
$arrA[58,0] = "20"
$arrB[58,0] = 0.0773
$arrC[58,0] = 0.145
$arrA[59,0] = "2"
$arrB[59,0] = 0.0263
$arrC[59,0] = 0.0951
$arrA[60,0] = "19"
$arrB[60,0] = 0.0707
$arrC[60,0] = 0.1699
$arrA[61,0] = "17"
$arrB[61,0] = 0.022
$arrC[61,0] = 0.3866
$arrA[62,0] = "16"
$arrB[62,0] = 0.1082
$arrC[62,0] = 0.6686
$arrA[63,0] = "14"
$arrB[63,0] = 0.039
$arrC[63,0] = 0.6963
$arrA[64,0] = "13"
$arrB[64,0] = 0.0714
$arrC[64,0] = 0.4161
$arrA[65,0] = "12"
$arrB[65,0] = 0.4584
$arrC[65,0] = 0.3191
$arrA[66,0] = "11"
$arrB[66,0] = -0.349
$arrC[66,0] = 0.6813
$arrA[67,0] = "10"
$arrB[67,0] = -0.1548
$arrC[67,0] = 1.2659
$arrA[68,0] = "1"
$arrB[68,0] = 0.0023
$arrC[68,0] = 0.0964
$arrA[69,0] = "350"
$arrB[69,0] = 0.0718
$arrC[69,0] = 1.023
$arrA[70,0] = "348"
$arrB[70,0] = 0.0008
$arrC[70,0] = 0.1673
$arrA[71,0] = "347"
$arrB[71,0] = 0.0229
$arrC[71,0] = 0.3054
$arrA[72,0] = "346"
$arrB[72,0] = 0.0023
$arrC[72,0] = 0.0749
$arrA[73,0] = "345"
$arrB[73,0] = 0.0827
$arrC[73,0] = 0.1248
$arrA[74,0] = "343"
$arrB[74,0] = -0.184
$arrC[74,0] = 0.296
$arrA[75,0] = "342"
$arrB[75,0] = 0.0843
$arrC[75,0] = 0.2048
$arrA[76,0] = "341"
$arrB[76,0] = 0.0591
$arrC[76,0] = 0.2658
$arrA[77,0] = "340"
$arrB[77,0] = 0.2156
$arrC[77,0] = 0.5736
$arrA[78,0] = "338"
$arrB[78,0] = -0.086
$arrC[78,0] = 0.6538
$arrA[79,0] = "337"
$arrB[79,0] = -0.0892
$arrC[79,0] = 0.5678
$arrA[80,0] = "335"
$arrB[80,0] = -0.0172
$arrC[80,0] = 0.2505
$arrA[81,0] = "334"
$arrB[81,0] = -0.0036
$arrC[81,0] = 0.2124
$arrA[82,0] = "333"
$arrB[82,0] = 0.0046
$arrC[82,0] = 0.4445
$arrA[83,0] = "332"
$arrB[83,0] = 0.0574
$arrC[83,0] = 0.1032
$arrA[84,0] = "331"
$arrB[84,0] = 0.0321
$arrC[84,0] = 0.4845
$arrA[85,0] = "330"
$arrB[85,0] = 0.0131
$arrC[85,0] = 0.2559
$arrA[86,0] = "329"
$arrB[86,0] = 0.0326
$arrC[86,0] = 0.162
$arrA[87,0] = "328"
$arrB[87,0] = 0.0243
$arrC[87,0] = 0.2076
$arrA[88,0] = "326"
$arrB[88,0] = -0.014
$arrC[88,0] = 0.0502
$arrA[89,0] = "325"
$arrB[89,0] = 0.1078
$arrC[89,0] = 0.2098
$arrA[90,0] = "324"
$arrB[90,0] = 0.043
$arrC[90,0] = 0.3091
$arrA[91,0] = "323"
$arrB[91,0] = 0.0165
$arrC[91,0] = 0.1048
$arrA[92,0] = "320"
$arrB[92,0] = 0.017
$arrC[92,0] = 0.2305
$arrA[93,0] = "319"
$arrB[93,0] = -0.0162
$arrC[93,0] = 0.1248
$arrA[94,0] = "317"
$arrB[94,0] = 0.0746
$arrC[94,0] = 0.5271
$arrA[95,0] = "316"
$arrB[95,0] = -0.0811
$arrC[95,0] = 0.6222
$arrA[96,0] = "315"
$arrB[96,0] = -0.0734
$arrC[96,0] = 0.6687
$arrA[97,0] = "314"
$arrB[97,0] = -0.0257
$arrC[97,0] = 0.3176
$arrA[98,0] = "313"
$arrB[98,0] = -0.0505
$arrC[98,0] = 0.1028
$arrA[99,0] = "312"
$arrB[99,0] = 0.0003
$arrC[99,0] = 0.2077
$arrA[100,0] = "308"
$arrB[100,0] = 0.0013
$arrC[100,0] = 0.543
$arrA[101,0] = "307"
$arrB[101,0] = 0.1328
$arrC[101,0] = 0.2289
$arrA[102,0] = "306"
$arrB[102,0] = 0.0168
$arrC[102,0] = 0.1262
$arrA[103,0] = "305"
$arrB[103,0] = -0.005
$arrC[103,0] = 0.0974
$arrA[104,0] = "304"
$arrB[104,0] = -0.0273
$arrC[104,0] = 0.0964
$arrA[105,0] = "303"
$arrB[105,0] = -0.0144
$arrC[105,0] = 0.142
$arrA[106,0] = "302"
$arrB[106,0] = 0.0236
$arrC[106,0] = 0.1607
$arrA[107,0] = "98"
$arrB[107,0] = -0.0134
$arrC[107,0] = 0.1563
$arrA[108,0] = "97"
$arrB[108,0] = 0.0622
$arrC[108,0] = 0.6319
$arrA[109,0] = "96"
$arrB[109,0] = -0.0101
$arrC[109,0] = 0.2069
$arrA[110,0] = "95"
$arrB[110,0] = -0.0298
$arrC[110,0] = 0.354
$arrA[111,0] = "92"
$arrB[111,0] = 0.0096
$arrC[111,0] = 0.103
$arrA[112,0] = "91"
$arrB[112,0] = -0.1687
$arrC[112,0] = 0.387
$arrA[113,0] = "90"
$arrB[113,0] = -0.0593
$arrC[113,0] = 0.1746
$arrA[114,0] = "89"
$arrB[114,0] = -0.0218
$arrC[114,0] = 0.088
$arrA[115,0] = "87"
$arrB[115,0] = 0.0142
$arrC[115,0] = 0.8121
$arrA[116,0] = "83"
$arrB[116,0] = 0.0169
$arrC[116,0] = 0.1131
$arrA[117,0] = "82"
$arrB[117,0] = 0.0519
$arrC[117,0] = 0.2238
$arrA[118,0] = "81"
$arrB[118,0] = -0.0116
$arrC[118,0] = 0.0756
$arrA[119,0] = "80"
$arrB[119,0] = 0.29
$arrC[119,0] = 0.5531
$arrA[120,0] = "79"
$arrB[120,0] = -0.0112
$arrC[120,0] = 0.4234
$arrA[121,0] = "78"
$arrB[121,0] = 0.0774
$arrC[121,0] = 0.0999
$arrA[122,0] = "75"
$arrB[122,0] = -0.0238
$arrC[122,0] = 0.1022
$arrA[123,0] = "74"
$arrB[123,0] = 0.0105
$arrC[123,0] = 0.097
$arrA[124,0] = "73"
$arrB[124,0] = -0.0303
$arrC[124,0] = 0.2122
$arrA[125,0] = "71"
$arrB[125,0] = -0.0071
$arrC[125,0] = 0.1928
$arrA[126,0] = "70"
$arrB[126,0] = -0.0006
$arrC[126,0] = 0.281
$arrA[127,0] = "69"
$arrB[127,0] = -0.0116
$arrC[127,0] = 0.1011
$arrA[128,0] = "67"
$arrB[128,0] = 0.0799
$arrC[128,0] = 0.2055
$arrA[129,0] = "66"
$arrB[129,0] = 0.1967
$arrC[129,0] = 0.3307
$arrA[130,0] = "65"
$arrB[130,0] = 0.2237
$arrC[130,0] = 0.4687
$arrA[131,0] = "64"
$arrB[131,0] = 0.0245
$arrC[131,0] = 0.1917
$arrA[132,0] = "63"
$arrB[132,0] = -0.0304
$arrC[132,0] = 0.204
$arrA[133,0] = "62"
$arrB[133,0] = -0.016
$arrC[133,0] = 0.1313
$arrA[134,0] = "56"
$arrB[134,0] = -0.0202
$arrC[134,0] = 0.1259
$arrA[135,0] = "55"
$arrB[135,0] = -0.0022
$arrC[135,0] = 0.1559
$arrA[136,0] = "54"
$arrB[136,0] = -0.0159
$arrC[136,0] = 0.3088
$arrA[137,0] = "53"
$arrB[137,0] = 0.0043
$arrC[137,0] = 0.1766
$arrA[138,0] = "52"
$arrB[138,0] = 0.1254
$arrC[138,0] = 0.129
$arrA[139,0] = "51"
$arrB[139,0] = 0.106
$arrC[139,0] = 0.397
$arrA[140,0] = "100"
$arrB[140,0] = -0.5258
$arrC[140,0] = 0.6429
$arrA[141,0] = "149"
$arrB[141,0] = 0.0041
$arrC[141,0] = 0.1015
$arrA[142,0] = "148"
$arrB[142,0] = 0.0455
$arrC[142,0] = 0.1755
$arrA[143,0] = "147"
$arrB[143,0] = 0.0085
$arrC[143,0] = 0.25
$arrA[144,0] = "146"
$arrB[144,0] = -0.0256
$arrC[144,0] = 0.281
$arrA[145,0] = "145"
$arrB[145,0] = 0.0421
$arrC[145,0] = 0.1244
$arrA[146,0] = "143"
$arrB[146,0] = -0.0091
$arrC[146,0] = 0.1053
$arrA[147,0] = "142"
$arrB[147,0] = 0.0309
$arrC[147,0] = 0.4025
$arrA[148,0] = "141"
$arrB[148,0] = 0.021
$arrC[148,0] = 0.1019
$arrA[149,0] = "140"
$arrB[149,0] = -0.0882
$arrC[149,0] = 0.2145
$arrA[150,0] = "139"
$arrB[150,0] = 0.0397
$arrC[150,0] = 0.183
$arrA[151,0] = "138"
$arrB[151,0] = 0.0654
$arrC[151,0] = 0.326
$arrA[152,0] = "137"
$arrB[152,0] = -0.0016
$arrC[152,0] = 0.1722
$arrA[153,0] = "136"
$arrB[153,0] = -0.0136
$arrC[153,0] = 0.2111
$arrA[154,0] = "135"
$arrB[154,0] = 0.0028
$arrC[154,0] = 0.0764
$arrA[155,0] = "134"
$arrB[155,0] = -0.018
$arrC[155,0] = 0.1091
$arrA[156,0] = "133"
$arrB[156,0] = 0.0405
$arrC[156,0] = 0.181
$arrA[157,0] = "132"
$arrB[157,0] = 0.053
$arrC[157,0] = 0.1515
$arrA[158,0] = "131"
$arrB[158,0] = -0.012
$arrC[158,0] = 0.2976
$arrA[159,0] = "130"
$arrB[159,0] = -0.0181
$arrC[159,0] = 0.1215
$arrA[160,0] = "129"
$arrB[160,0] = -0.0335
$arrC[160,0] = 0.2029
$arrA[161,0] = "128"
$arrB[161,0] = 0.0014
$arrC[161,0] = 0.335
$arrA[162,0] = "125"
$arrB[162,0] = -0.0154
$arrC[162,0] = 0.2669
$arrA[163,0] = "124"
$arrB[163,0] = 0.0055
$arrC[163,0] = 0.1817
$arrA[164,0] = "122"
$arrB[164,0] = -0.0004
$arrC[164,0] = 0.1213
$arrA[165,0] = "121"
$arrB[165,0] = -0.3033
$arrC[165,0] = 0.1308
$arrA[166,0] = "120"
$arrB[166,0] = -0.0887
$arrC[166,0] = 0.4833
$arrA[167,0] = "117"
$arrB[167,0] = 0.2813
$arrC[167,0] = 0.3482
$arrA[168,0] = "116"
$arrB[168,0] = -0.0106
$arrC[168,0] = 0.2364
$arrA[169,0] = "115"
$arrB[169,0] = -0.0675
$arrC[169,0] = 0.3253
$arrA[170,0] = "114"
$arrB[170,0] = -0.039
$arrC[170,0] = 0.2977
$arrA[171,0] = "111"
$arrB[171,0] = -0.0144
$arrC[171,0] = 0.0658
$arrA[172,0] = "108"
$arrB[172,0] = -0.0313
$arrC[172,0] = 0.3863
$arrA[173,0] = "107"
$arrB[173,0] = -0.0893
$arrC[173,0] = 0.7023
$arrA[174,0] = "106"
$arrB[174,0] = 0.0169
$arrC[174,0] = 0.7485
$arrA[175,0] = "104"
$arrB[175,0] = 0.0008
$arrC[175,0] = 0.4123
$arrA[176,0] = "103"
$arrB[176,0] = -0.0085
$arrC[176,0] = 0.1279
$arrA[177,0] = "101"
$arrB[177,0] = -0.4228
$arrC[177,0] = 0.3811
$arrA[178,0] = "383"
$arrB[178,0] = -0.0841
$arrC[178,0] = 0.189
$arrA[179,0] = "382"
$arrB[179,0] = 0.0934
$arrC[179,0] = 0.1577
$arrA[180,0] = "381"
$arrB[180,0] = 0.092
$arrC[180,0] = 0.17
$arrA[181,0] = "380"
$arrB[181,0] = 0.0241
$arrC[181,0] = 0.3829
$arrA[182,0] = "379"
$arrB[182,0] = 0.0337
$arrC[182,0] = 0.1055
$arrA[183,0] = "378"
$arrB[183,0] = 0.0251
$arrC[183,0] = 0.3386
$arrA[184,0] = "377"
$arrB[184,0] = -0.0052
$arrC[184,0] = 0.326
$arrA[185,0] = "376"
$arrB[185,0] = 0.0234
$arrC[185,0] = 0.1906
$arrA[186,0] = "374"
$arrB[186,0] = 0.0561
$arrC[186,0] = 0.2073
$arrA[187,0] = "373"
$arrB[187,0] = 0.01
$arrC[187,0] = 0.2381
$arrA[188,0] = "372"
$arrB[188,0] = 0.057
$arrC[188,0] = 0.1793
$arrA[189,0] = "370"
$arrB[189,0] = 0.0057
$arrC[189,0] = 0.1562
$arrA[190,0] = "369"
$arrB[190,0] = 0.0102
$arrC[190,0] = 0.135
$arrA[191,0] = "367"
$arrB[191,0] = 0.0043
$arrC[191,0] = 0.215
$arrA[192,0] = "364"
$arrB[192,0] = -0.0006
$arrC[192,0] = 0.0891
$arrA[193,0] = "362"
$arrB[193,0] = -0.1235
$arrC[193,0] = 0.3601
$arrA[194,0] = "361"
$arrB[194,0] = 0.0111
$arrC[194,0] = 0.5265
$arrA[195,0] = "360"
$arrB[195,0] = -0.0333
$arrC[195,0] = 0.2992
$arrA[196,0] = "359"
$arrB[196,0] = -0.058
$arrC[196,0] = 0.2198
$arrA[197,0] = "354"
$arrB[197,0] = -0.0299
$arrC[197,0] = 0.1044
$arrA[198,0] = "353"
$arrB[198,0] = 0.0156
$arrC[198,0] = 0.4178
$arrA[199,0] = "352"
$arrB[199,0] = 0.0205
$arrC[199,0] = 0.3503
$arrA[200,0] = "351"
$arrB[200,0] = -0.0818
$arrC[200,0] = 0.1863
$arrA[201,0] = "196"
$arrB[201,0] = 0.2331
$arrC[201,0] = 0.803
$arrA[202,0] = "195"
$arrB[202,0] = -0.0001
$arrC[202,0] = 0.5624
$arrA[203,0] = "192"
$arrB[203,0] = 0.0723
$arrC[203,0] = 0.2842
$arrA[204,0] = "190"
$arrB[204,0] = -0.039
$arrC[204,0] = 0.2112
$arrA[205,0] = "189"
$arrB[205,0] = 0.0823
$arrC[205,0] = 1.0237
$arrA[206,0] = "188"
$arrB[206,0] = -0.0226
$arrC[206,0] = 0.3139
$arrA[207,0] = "187"
$arrB[207,0] = 0.0373
$arrC[207,0] = 0.7147
$arrA[208,0] = "183"
$arrB[208,0] = -0.0356
$arrC[208,0] = 0.2611
$arrA[209,0] = "182"
$arrB[209,0] = 0.0122
$arrC[209,0] = 0.4192
$arrA[210,0] = "181"
$arrB[210,0] = -0.0067
$arrC[210,0] = 0.1714
$arrA[211,0] = "179"
$arrB[211,0] = -0.013
$arrC[211,0] = 0.2407
$arrA[212,0] = "178"
$arrB[212,0] = 0.1062
$arrC[212,0] = 0.2668
$arrA[213,0] = "175"
$arrB[213,0] = -0.0123
$arrC[213,0] = 0.1031
$arrA[214,0] = "174"
$arrB[214,0] = 0.0334
$arrC[214,0] = 0.4104
$arrA[215,0] = "172"
$arrB[215,0] = 0.0259
$arrC[215,0] = 0.5338
$arrA[216,0] = "171"
$arrB[216,0] = 0.0339
$arrC[216,0] = 0.2011
$arrA[217,0] = "169"
$arrB[217,0] = 0.0074
$arrC[217,0] = 0.2799
$arrA[218,0] = "168"
$arrB[218,0] = -0.0488
$arrC[218,0] = 0.4896
$arrA[219,0] = "166"
$arrB[219,0] = 0.0932
$arrC[219,0] = 0.4142
$arrA[220,0] = "165"
$arrB[220,0] = -0.0651
$arrC[220,0] = 0.251
$arrA[221,0] = "164"
$arrB[221,0] = -0.0039
$arrC[221,0] = 0.1631
$arrA[222,0] = "163"
$arrB[222,0] = -0.0002
$arrC[222,0] = 0.5504
$arrA[223,0] = "161"
$arrB[223,0] = 0.0716
$arrC[223,0] = 0.261
$arrA[224,0] = "160"
$arrB[224,0] = -0.0101
$arrC[224,0] = 0.0964
$arrA[225,0] = "158"
$arrB[225,0] = -0.0363
$arrC[225,0] = 0.2019
$arrA[226,0] = "157"
$arrB[226,0] = -0.023
$arrC[226,0] = 0.2126
$arrA[227,0] = "156"
$arrB[227,0] = 0.0742
$arrC[227,0] = 0.1728
$arrA[228,0] = "155"
$arrB[228,0] = 0.0032
$arrC[228,0] = 0.1725
$arrA[229,0] = "154"
$arrB[229,0] = 0.0014
$arrC[229,0] = 0.1547
$arrA[230,0] = "153"
$arrB[230,0] = -0.1379
$arrC[230,0] = 0.5379
$arrA[231,0] = "152"
$arrB[231,0] = 0.1193
$arrC[231,0] = 0.331
$arrA[232,0] = "247"
$arrB[232,0] = 0.0259
$arrC[232,0] = 0.1323
$arrA[233,0] = "245"
$arrB[233,0] = 0.2104
$arrC[233,0] = 0.1408
$arrA[234,0] = "243"
$arrB[234,0] = 0.0729
$arrC[234,0] = 0.236
$arrA[235,0] = "241"
$arrB[235,0] = 0.0148
$arrC[235,0] = 0.3174
$arrA[236,0] = "240"
$arrB[236,0] = -0.077
$arrC[236,0] = 0.3698
$arrA[237,0] = "239"
$arrB[237,0] = 0.0223
$arrC[237,0] = 0.4053
$arrA[238,0] = "238"
$arrB[238,0] = 0.0724
$arrC[238,0] = 0.423
$arrA[239,0] = "237"
$arrB[239,0] = -0.0548
$arrC[239,0] = 0.287
$arrA[240,0] = "234"
$arrB[240,0] = -0.0317
$arrC[240,0] = 0.0748
$arrA[241,0] = "233"
$arrB[241,0] = -0.0402
$arrC[241,0] = 0.1228
$arrA[242,0] = "231"
$arrB[242,0] = -0.3102
$arrC[242,0] = 0.434
$arrA[243,0] = "230"
$arrB[243,0] = 0.0046
$arrC[243,0] = 0.4236
$arrA[244,0] = "229"
$arrB[244,0] = -0.0697
$arrC[244,0] = 0.1841
$arrA[245,0] = "228"
$arrB[245,0] = 0.0239
$arrC[245,0] = 0.3507
$arrA[246,0] = "226"
$arrB[246,0] = -0.0707
$arrC[246,0] = 0.1214
$arrA[247,0] = "225"
$arrB[247,0] = -0.3359
$arrC[247,0] = 0.3753
$arrA[248,0] = "224"
$arrB[248,0] = -0.0297
$arrC[248,0] = 0.1816
$arrA[249,0] = "223"
$arrB[249,0] = 0.089
$arrC[249,0] = 0.0594
$arrA[250,0] = "222"
$arrB[250,0] = 0.0855
$arrC[250,0] = 0.5541
$arrA[251,0] = "220"
$arrB[251,0] = 0.024
$arrC[251,0] = 0.1259
$arrA[252,0] = "219"
$arrB[252,0] = -0.1084
$arrC[252,0] = 0.2093
$arrA[253,0] = "217"
$arrB[253,0] = 0.0812
$arrC[253,0] = 0.8228
$arrA[254,0] = "216"
$arrB[254,0] = -0.0236
$arrC[254,0] = 0.1392
$arrA[255,0] = "215"
$arrB[255,0] = 0.0056
$arrC[255,0] = 0.3164
$arrA[256,0] = "214"
$arrB[256,0] = -0.2921
$arrC[256,0] = 0.7268
$arrA[257,0] = "213"
$arrB[257,0] = 0.0625
$arrC[257,0] = 0.578
$arrA[258,0] = "212"
$arrB[258,0] = -0.0058
$arrC[258,0] = 0.1894
$arrA[259,0] = "211"
$arrB[259,0] = 0.0246
$arrC[259,0] = 0.407
$arrA[260,0] = "210"
$arrB[260,0] = -0.1849
$arrC[260,0] = 0.4927
$arrA[261,0] = "209"
$arrB[261,0] = -0.1914
$arrC[261,0] = 0.4267
$arrA[262,0] = "208"
$arrB[262,0] = -0.0163
$arrC[262,0] = 0.2164
$arrA[263,0] = "207"
$arrB[263,0] = 0.0085
$arrC[263,0] = 0.1103
$arrA[264,0] = "205"
$arrB[264,0] = 0.0374
$arrC[264,0] = 0.2088
$arrA[265,0] = "202"
$arrB[265,0] = 0.0031
$arrC[265,0] = 0.1658
$arrA[266,0] = "201"
$arrB[266,0] = -0.04
$arrC[266,0] = 0.162
$arrA[267,0] = "400"
$arrB[267,0] = 0.0136
$arrC[267,0] = 0.1265
$arrA[268,0] = "399"
$arrB[268,0] = 0.0556
$arrC[268,0] = 0.3509
$arrA[269,0] = "398"
$arrB[269,0] = -0.0126
$arrC[269,0] = 0.1297
$arrA[270,0] = "396"
$arrB[270,0] = -0.0135
$arrC[270,0] = 0.0814
$arrA[271,0] = "395"
$arrB[271,0] = 0.0053
$arrC[271,0] = 0.1715
$arrA[272,0] = "394"
$arrB[272,0] = 0.1002
$arrC[272,0] = 0.1516
$arrA[273,0] = "393"
$arrB[273,0] = 0.0943
$arrC[273,0] = 0.3232
$arrA[274,0] = "392"
$arrB[274,0] = 0.0029
$arrC[274,0] = 0.2025
$arrA[275,0] = "391"
$arrB[275,0] = 0.0088
$arrC[275,0] = 0.1519
$arrA[276,0] = "390"
$arrB[276,0] = 0.0275
$arrC[276,0] = 0.1196
$arrA[277,0] = "389"
$arrB[277,0] = -0.0043
$arrC[277,0] = 0.4332
$arrA[278,0] = "388"
$arrB[278,0] = 0.0393
$arrC[278,0] = 0.2056
$arrA[279,0] = "387"
$arrB[279,0] = 0.0325
$arrC[279,0] = 0.1707
$arrA[280,0] = "386"
$arrB[280,0] = -0.0033
$arrC[280,0] = 0.1014
$arrA[281,0] = "385"
$arrB[281,0] = -0.0164
$arrC[281,0] = 0.0509

$ws.Range("A2:A283").Value = $arrA
$ws.Range("B2:B283").Value = $arrB
$ws.Range("C2:C283").Value = $arrC